# Remove columns from processed data: delete whole columns B ("b") and D ("d")
# from the sheet. Deleting B first shifts D into C, so we then delete the
# (new) column C to remove the original "d" column. All remaining columns
# shift left accordingly, matching the target layout (A1:L12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).EntireColumn.Delete()
$ws.Columns.Item(3).EntireColumn.Delete()
